$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 8) with the same date/text as the most recent
# existing entry (26-09-2025), reusing the existing shared-string values.
$ws.Range("A8").Value = "26-09-2025"
$ws.Range("B8").Value = "The price of gold in India today is ₹11,488 per gram for 24 karat gold, ₹10,530 per gram for 22 karat gold and ₹8,616 per gram for 18 karat gold (also called 999 gold)."
